# The target diff for this document is purely a cosmetic re-serialization
# of the underlying OOXML: namespace declarations and element attributes
# were reordered (alphabetically) and no text, formatting, structure, or
# any other semantic content was added, removed, or modified anywhere in
# word/document.xml, word/footer1.xml, or word/styles.xml. Every removed
# line in the diff has an exact counterpart added line with the same tag
# name and the same set of attribute name/value pairs, just written in a
# different order (this matches the "Moving from 2.0.0 to 2.0.1" commit,
# i.e. a tooling/library version bump that mechanically re-saved the test
# fixture without touching its content).
#
# Word's COM object model does not expose any control over the low level
# attribute/namespace ordering used when OOXML parts are serialized, so
# there is no content-level edit to perform here: the correct action is
# to leave the document exactly as it is.

$d = $word.ActiveDocument
